$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.108.05"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "1.599.17"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'212.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").Value = "'17.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "'0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").Value = "1.821.68"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").Value = "1.594.33"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "'0.511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "26.078.34"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "'60.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "'204.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.36%  "
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("E24").Value = "  +11.66%  "
$ws.Range("D25").Value = "'141.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").Value = "'15.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("D29").Value = "'6.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "'1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "1.110.31"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("E37").Value = "  +8.17%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").Value = "'0.777"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "'0.780"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").Value = "1.733.10"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'92.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").Value = "'1.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").Value = "'53.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  +1.35%  "
